# Update pl_mw.xlsx (Case_2_20, 380 kV case) result values in Sheet1
# Rows 2-25 (data rows), columns B, C, E, F, G, I, L, N, O get new simulation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9992429780979819
$ws.Range("C2").Value = 0.1811929599846565
$ws.Range("E2").Value = 0.08941069587415384
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002415944870811899
$ws.Range("I2").Value = 0.5666943818440977
$ws.Range("L2").Value = 0.2183729955876572
$ws.Range("N2").Value = 1.092475649779793
$ws.Range("O2").Value = 2.199927829652808
$ws.Range("B3").Value = 0.9044574867328379
$ws.Range("C3").Value = 0.1701148365688567
$ws.Range("E3").Value = 0.08987622655237715
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002418578993068641
$ws.Range("I3").Value = 0.5730978903180102
$ws.Range("L3").Value = 0.2090267239320553
$ws.Range("N3").Value = 1.097425316707159
$ws.Range("O3").Value = 2.205771986625706
$ws.Range("B4").Value = 0.8463364081106306
$ws.Range("C4").Value = 0.1632587100100835
$ws.Range("E4").Value = 0.09020801821514013
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002420283332495748
$ws.Range("I4").Value = 0.5774106140121127
$ws.Range("L4").Value = 0.2033976418929626
$ws.Range("N4").Value = 1.100909681150497
$ws.Range("O4").Value = 2.210976707920167
$ws.Range("B5").Value = 0.8226726384100971
$ws.Range("C5").Value = 0.1604513270380608
$ws.Range("E5").Value = 0.0903547884998428
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002420999800698593
$ws.Range("I5").Value = 0.5792637569943047
$ws.Range("L5").Value = 0.2011313865170337
$ws.Range("N5").Value = 1.102441746076913
$ws.Range("O5").Value = 2.213503554030865
$ws.Range("B6").Value = 0.8187446034462482
$ws.Range("C6").Value = 0.1599843550161211
$ws.Range("E6").Value = 0.09037985814584104
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002421120096471237
$ws.Range("I6").Value = 0.5795772449274139
$ws.Range("L6").Value = 0.2007567486563318
$ws.Range("N6").Value = 1.102702924235615
$ws.Range("O6").Value = 2.213947632276998
$ws.Range("B7").Value = 0.8460171830114405
$ws.Range("C7").Value = 0.1632209029413048
$ws.Range("E7").Value = 0.09020995078493499
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002420292906208659
$ws.Range("I7").Value = 0.5774352189275014
$ws.Range("L7").Value = 0.203366966342287
$ws.Range("N7").Value = 1.100929888726959
$ws.Range("O7").Value = 2.211009143221574
$ws.Range("B8").Value = 0.966545722494061
$ws.Range("C8").Value = 0.1773845421386966
$ws.Range("E8").Value = 0.08956167720910457
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002416835103094652
$ws.Range("I8").Value = 0.5688231757028426
$ws.Range("L8").Value = 0.2151277111049552
$ws.Range("N8").Value = 1.094089999098429
$ws.Range("O8").Value = 2.201607063950888
$ws.Range("B9").Value = 1.203462094713529
$ws.Range("C9").Value = 0.2047245902438988
$ws.Range("E9").Value = 0.08865478370805135
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002410741522460305
$ws.Range("I9").Value = 0.554962981134107
$ws.Range("L9").Value = 0.2390573153380871
$ws.Range("N9").Value = 1.084201889371883
$ws.Range("O9").Value = 2.196023508471228
$ws.Range("B10").Value = 1.377808985267734
$ws.Range("C10").Value = 0.224540760417284
$ws.Range("E10").Value = 0.08821035337419936
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002406679321737776
$ws.Range("I10").Value = 0.5466333523415656
$ws.Range("L10").Value = 0.2571654938276851
$ws.Range("N10").Value = 1.07907626834843
$ws.Range("O10").Value = 2.199800887124383
$ws.Range("B11").Value = 1.457174637251399
$ws.Range("C11").Value = 0.2334958465907562
$ws.Range("E11").Value = 0.08805629602954568
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002404920496196962
$ws.Range("I11").Value = 0.543248076673116
$ws.Range("L11").Value = 0.2655176706438311
$ws.Range("N11").Value = 1.077207089805412
$ws.Range("O11").Value = 2.203239909998246
$ws.Range("B12").Value = 1.487234849216293
$ws.Range("C12").Value = 0.2368782346943306
$ws.Range("E12").Value = 0.08800487232127807
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002404267218942185
$ws.Range("I12").Value = 0.5420243791066817
$ws.Range("L12").Value = 0.2686968522414332
$ws.Range("N12").Value = 1.076565627742738
$ws.Range("O12").Value = 2.204790315628799
$ws.Range("B13").Value = 1.480760597012363
$ws.Range("C13").Value = 0.2361501668977439
$ws.Range("E13").Value = 0.08801563988771832
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002404407347672451
$ws.Range("I13").Value = 0.5422853319564638
$ws.Range("L13").Value = 0.2680114310385875
$ws.Range("N13").Value = 1.076700829251024
$ws.Range("O13").Value = 2.204445361743865
$ws.Range("B14").Value = 1.459647594964395
$ws.Range("C14").Value = 0.2337742931936759
$ws.Range("E14").Value = 0.0880519268202562
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002404866495365616
$ws.Range("I14").Value = 0.5431462344035012
$ws.Range("L14").Value = 0.2657788958319287
$ws.Range("N14").Value = 1.077152987321625
$ws.Range("O14").Value = 2.203362486110109
$ws.Range("B15").Value = 1.446716024940542
$ws.Range("C15").Value = 0.2323178628266191
$ws.Range("E15").Value = 0.08807505394488579
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002405149395747318
$ws.Range("I15").Value = 0.5436811501096912
$ws.Range("L15").Value = 0.2644135355980524
$ws.Range("N15").Value = 1.077438584227963
$ws.Range("O15").Value = 2.202731526845952
$ws.Range("B16").Value = 1.372623205699256
$ws.Range("C16").Value = 0.2239543149925112
$ws.Range("E16").Value = 0.08822138921400224
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002406796052641217
$ws.Range("I16").Value = 0.546862728724026
$ws.Range("L16").Value = 0.2566219587357637
$ws.Range("N16").Value = 1.079207717108133
$ws.Range("O16").Value = 2.19961082241673
$ws.Range("B17").Value = 1.327182407739429
$ws.Range("C17").Value = 0.2188082182093183
$ws.Range("E17").Value = 0.08832348162729708
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002407828999324967
$ws.Range("I17").Value = 0.5489180834249261
$ws.Range("L17").Value = 0.2518713811430615
$ws.Range("N17").Value = 1.08041136951735
$ws.Range("O17").Value = 2.198137570014808
$ws.Range("B18").Value = 1.301051279648163
$ws.Range("C18").Value = 0.2158427395233389
$ws.Range("E18").Value = 0.08838673168881606
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002408431511859671
$ws.Range("I18").Value = 0.5501382784223772
$ws.Range("L18").Value = 0.249149774169112
$ws.Range("N18").Value = 1.081147219181048
$ws.Range("O18").Value = 2.197452111012353
$ws.Range("B19").Value = 1.292204675268295
$ws.Range("C19").Value = 0.2148377261818553
$ws.Range("E19").Value = 0.08840892510757037
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002408636955028332
$ws.Range("I19").Value = 0.5505579390353645
$ws.Range("L19").Value = 0.248230142596654
$ws.Range("N19").Value = 1.081403848047614
$ws.Range("O19").Value = 2.197247814642253
$ws.Range("B20").Value = 1.332019130328547
$ws.Range("C20").Value = 0.2193566077225455
$ws.Range("E20").Value = 0.08831214499318207
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002407718172430406
$ws.Range("I20").Value = 0.5486953524929383
$ws.Range("L20").Value = 0.2523759712188678
$ws.Range("N20").Value = 1.080278733475794
$ws.Range("O20").Value = 2.198277637445671
$ws.Range("B21").Value = 1.465848841406
$ws.Range("C21").Value = 0.2344723825860626
$ws.Range("E21").Value = 0.08804108084879481
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002404731286985283
$ws.Range("I21").Value = 0.5428917848649419
$ws.Range("L21").Value = 0.2664342013968621
$ws.Range("N21").Value = 1.077018377847878
$ws.Range("O21").Value = 2.203673813465826
$ws.Range("B22").Value = 1.553349444657783
$ws.Range("C22").Value = 0.2443005765226474
$ws.Range("E22").Value = 0.08790422540030107
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.00240285348500561
$ws.Range("I22").Value = 0.5394383132673681
$ws.Range("L22").Value = 0.2757175763149178
$ws.Range("N22").Value = 1.07527424992297
$ws.Range("O22").Value = 2.208647092347491
$ws.Range("B23").Value = 1.506646056237173
$ws.Range("C23").Value = 0.2390597917822959
$ws.Range("E23").Value = 0.0879735817423537
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002403848924169439
$ws.Range("I23").Value = 0.5412503825466821
$ws.Range("L23").Value = 0.2707541573640242
$ws.Range("N23").Value = 1.076169788957344
$ws.Range("O23").Value = 2.205860171083543
$ws.Range("B24").Value = 1.329832468502502
$ws.Range("C24").Value = 0.2191087023642808
$ws.Range("E24").Value = 0.0883172560946619
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.0024077682503356
$ws.Range("I24").Value = 0.5487959290949718
$ws.Range("L24").Value = 0.2521478162507123
$ws.Range("N24").Value = 1.080338561615093
$ws.Range("O24").Value = 2.198213809842969
$ws.Range("B25").Value = 1.139315879153003
$ws.Range("C25").Value = 0.1973755015560243
$ws.Range("E25").Value = 0.08886114153467339
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002412316863941194
$ws.Range("I25").Value = 0.5583875868166039
$ws.Range("L25").Value = 0.2324910620599212
$ws.Range("N25").Value = 1.08650054285652
$ws.Range("O25").Value = 2.196023508471228
